$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at X (column 24). This shifts the existing X (nom)
# and Y (url_produit) columns one position to the right (to Y and Z),
# and extends the used range from A1:Y205 to A1:Z205.
$ws.Range("X1").EntireColumn.Insert()

# New header for the inserted column: a timestamp, consistent with the
# other snapshot-date headers in row 1 (B1..W1).
$ws.Cells.Item(1, 24).Value = "2026-01-28 17:24:56"

# For each data row, the new X column should repeat the latest price
# already recorded in column W (the previous last snapshot column) --
# but only for rows where a price exists (rows 2-80). Rows 81-205 have
# no price data in W, so the new X cell is left blank there too.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $priceCell = $ws.Cells.Item($r, 23)
    $price = $priceCell.Value2()
    if ($price -ne $null -and $price -ne "") {
        $ws.Cells.Item($r, 24).Value2 = $price
    }
}
